$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 395.7073
$ws.Range("J17").Value = 403.1
$ws.Range("L17").Value = 1209.3
$ws.Range("N17").Value = -1545.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 703.5
$ws.Range("I33").Value = 629.7273
$ws.Range("J33").Value = 974
$ws.Range("K33").Value = 629.7273
$ws.Range("L33").Value = 974
$ws.Range("M33").Value = -400.7273
$ws.Range("N33").Value = -1432

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 492.25
$ws.Range("J41").Value = 439.66666
$ws.Range("L41").Value = 439.66666
$ws.Range("N41").Value = -1319.66666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5835.75
$ws.Range("I43").Value = 1400
$ws.Range("J43").Value = 6469.4287
$ws.Range("K43").Value = 1400
$ws.Range("L43").Value = 6469.4287
$ws.Range("M43").Value = -1331
$ws.Range("N43").Value = -6607.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 47623050
$ws.Range("I64").Value = 66670344
$ws.Range("K64").Value = 66670344
$ws.Range("M64").Value = -66670096

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 47623050
$ws.Range("I67").Value = 66670344
$ws.Range("K67").Value = 66670344
$ws.Range("M67").Value = -66669486

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 71430090
$ws.Range("I86").Value = 125001350
$ws.Range("J86").Value = 1749.1666
$ws.Range("K86").Value = 125001350
$ws.Range("L86").Value = 1749.1666
$ws.Range("M86").Value = -125000227
$ws.Range("N86").Value = -3995.1666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 71430090
$ws.Range("I89").Value = 125001350
$ws.Range("J89").Value = 1749.1666
$ws.Range("K89").Value = 625006750
$ws.Range("L89").Value = 8745.833000000001
$ws.Range("M89").Value = -625001134
$ws.Range("N89").Value = -19977.833

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2825.9285
$ws.Range("I137").Value = 2021.1428
$ws.Range("J137").Value = 5240.2856
$ws.Range("K137").Value = 6063.428400000001
$ws.Range("L137").Value = 15720.8568
$ws.Range("M137").Value = -3513.428400000001
$ws.Range("N137").Value = -20820.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4926.355
$ws.Range("I2").Value = 4939.3335
$ws.Range("K2").Value = 4939.3335
$ws.Range("M2").Value = -4826.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3881.1875
$ws.Range("I61").Value = 2499.6667
$ws.Range("K61").Value = 2499.6667
$ws.Range("M61").Value = -2287.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3791.875
$ws.Range("J63").Value = 5052.1
$ws.Range("L63").Value = 5052.1
$ws.Range("N63").Value = -6424.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3791.875
$ws.Range("J66").Value = 5052.1
$ws.Range("L66").Value = 25260.5
$ws.Range("N66").Value = -32124.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 123408.164
$ws.Range("J92").Value = 123408.164
$ws.Range("L92").Value = 123408.164
$ws.Range("N92").Value = -128400.164

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 4926.355
$ws.Range("I116").Value = 4939.3335
$ws.Range("K116").Value = 4939.3335
$ws.Range("M116").Value = -2645.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1420494
$ws.Range("I132").Value = 2166636
$ws.Range("K132").Value = 6499908
$ws.Range("M132").Value = -6497378

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3881.1875
$ws.Range("I136").Value = 2499.6667
$ws.Range("K136").Value = 7499.000100000001
$ws.Range("M136").Value = -4949.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4926.355
$ws.Range("I3").Value = 4939.3335
$ws.Range("K3").Value = 4939.3335
$ws.Range("M3").Value = -4825.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 49999
$ws.Range("J92").Value = 49999
$ws.Range("L92").Value = 49999
$ws.Range("N92").Value = -54991

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7476.9434
$ws.Range("I99").Value = 7519.675
$ws.Range("K99").Value = 7519.675
$ws.Range("M99").Value = -6021.675

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7146934
$ws.Range("I107").Value = 12502760
$ws.Range("K107").Value = 12502760
$ws.Range("M107").Value = -12500840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 38467216
$ws.Range("I16").Value = 100003680
$ws.Range("J16").Value = 6924.625
$ws.Range("K16").Value = 100003680
$ws.Range("L16").Value = 6924.625
$ws.Range("M16").Value = -100003393
$ws.Range("N16").Value = -7498.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2624.9375
$ws.Range("J31").Value = 3037.4167
$ws.Range("L31").Value = 3037.4167
$ws.Range("N31").Value = -3627.4167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2624.9375
$ws.Range("J34").Value = 3037.4167
$ws.Range("L34").Value = 3037.4167
$ws.Range("N34").Value = -3441.4167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 38467216
$ws.Range("I113").Value = 100003680
$ws.Range("J113").Value = 6924.625
$ws.Range("K113").Value = 100003680
$ws.Range("L113").Value = 6924.625
$ws.Range("M113").Value = -100001510
$ws.Range("N113").Value = -11264.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4062.8
$ws.Range("I122").Value = 4032.5
$ws.Range("J122").Value = 4083
$ws.Range("K122").Value = 12097.5
$ws.Range("L122").Value = 12249
$ws.Range("M122").Value = -9647.5
$ws.Range("N122").Value = -17149

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5456.943
$ws.Range("I132").Value = 4638.4614
$ws.Range("J132").Value = 7821.4443
$ws.Range("K132").Value = 13915.3842
$ws.Range("L132").Value = 23464.3329
$ws.Range("M132").Value = -11385.3842
$ws.Range("N132").Value = -28524.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 32267534
$ws.Range("I134").Value = 125004730
$ws.Range("J134").Value = 11119
$ws.Range("K134").Value = 375014190
$ws.Range("L134").Value = 33357
$ws.Range("M134").Value = -375011655
$ws.Range("N134").Value = -38427

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 51460
$ws.Range("J9").Value = 58597.145
$ws.Range("L9").Value = 175791.435
$ws.Range("N9").Value = -176239.435

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 25
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 75
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 98
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 219919.7
$ws.Range("J37").Value = 219919.7
$ws.Range("L37").Value = 659759.1000000001
$ws.Range("N37").Value = -659983.1000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 71.25
$ws.Range("I38").Value = 35
$ws.Range("J38").Value = 83.333336
$ws.Range("K38").Value = 105
$ws.Range("L38").Value = 250.000008
$ws.Range("M38").Value = 242
$ws.Range("N38").Value = -944.000008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3472.7273
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 3720
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 11160
$ws.Range("M39").Value = -2706
$ws.Range("N39").Value = -11748

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 9261.666999999999
$ws.Range("I40").Value = 151.66667
$ws.Range("J40").Value = 13816.667
$ws.Range("K40").Value = 606.66668
$ws.Range("L40").Value = 55266.668
$ws.Range("M40").Value = -537.66668
$ws.Range("N40").Value = -55404.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1370.5834
$ws.Range("I98").Value = 3000.5
$ws.Range("J98").Value = 1044.6
$ws.Range("K98").Value = 9001.5
$ws.Range("L98").Value = 3133.8
$ws.Range("M98").Value = -7503.5
$ws.Range("N98").Value = -6129.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3097.3076
$ws.Range("I132").Value = 999
$ws.Range("J132").Value = 3478.818
$ws.Range("K132").Value = 8991
$ws.Range("L132").Value = 31309.362
$ws.Range("M132").Value = -6461
$ws.Range("N132").Value = -36369.362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10773
$ws.Range("I70").Value = 9064
$ws.Range("K70").Value = 9064
$ws.Range("M70").Value = -8794

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10773
$ws.Range("I73").Value = 9064
$ws.Range("K73").Value = 9064
$ws.Range("M73").Value = -8128

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 111119000
$ws.Range("I80").Value = 333335000
$ws.Range("J80").Value = 11000
$ws.Range("K80").Value = 333335000
$ws.Range("L80").Value = 11000
$ws.Range("M80").Value = -333334002
$ws.Range("N80").Value = -12996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 111119000
$ws.Range("I83").Value = 333335000
$ws.Range("J83").Value = 11000
$ws.Range("K83").Value = 1666675000
$ws.Range("L83").Value = 55000
$ws.Range("M83").Value = -1666670008
$ws.Range("N83").Value = -64984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6243.2104
$ws.Range("I113").Value = 2078.5
$ws.Range("J113").Value = 9272.091
$ws.Range("K113").Value = 2078.5
$ws.Range("L113").Value = 9272.091
$ws.Range("M113").Value = 91.5
$ws.Range("N113").Value = -13612.091

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 55000
$ws.Range("J127").Value = 55000
$ws.Range("L127").Value = 55000
$ws.Range("N127").Value = -64920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 40002892
$ws.Range("I132").Value = 58825870
$ws.Range("J132").Value = 4061.375
$ws.Range("K132").Value = 176477610
$ws.Range("L132").Value = 12184.125
$ws.Range("M132").Value = -176475080
$ws.Range("N132").Value = -17244.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1305.6154
$ws.Range("I93").Value = 1734.375
$ws.Range("J93").Value = 619.6
$ws.Range("K93").Value = 1734.375
$ws.Range("L93").Value = 619.6
$ws.Range("M93").Value = -486.375
$ws.Range("N93").Value = -3115.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 21546.5
$ws.Range("J54").Value = 22062
$ws.Range("L54").Value = 22062
$ws.Range("N54").Value = -23102

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2768.75
$ws.Range("I96").Value = 1744.3334
$ws.Range("J96").Value = 3005.1538
$ws.Range("K96").Value = 1744.3334
$ws.Range("L96").Value = 3005.1538
$ws.Range("M96").Value = -371.3334
$ws.Range("N96").Value = -5751.1538
